# Auto-generated edit script: updates LeveProfit calculation-dependent
# columns (H, I/J currentAveragePrice*, K/L LevePrice*, M/N LeveProfit*)
# across all 8 job sheets per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 819.5
$ws.Range("J28").Value = 926.6667
$ws.Range("L28").Value = 926.6667
$ws.Range("N28").Value = -1896.6667
$ws.Range("H40").Value = 1879.3478
$ws.Range("I40").Value = 1866.6666
$ws.Range("J40").Value = 1925
$ws.Range("K40").Value = 1866.6666
$ws.Range("L40").Value = 1925
$ws.Range("M40").Value = -1691.6666
$ws.Range("N40").Value = -2275
$ws.Range("H43").Value = 1019.1739
$ws.Range("J43").Value = 831.7646999999999
$ws.Range("L43").Value = 831.7646999999999
$ws.Range("N43").Value = -969.7646999999999
$ws.Range("H98").Value = 1293.0416
$ws.Range("I98").Value = 920.619
$ws.Range("J98").Value = 3900
$ws.Range("K98").Value = 920.619
$ws.Range("L98").Value = 3900
$ws.Range("M98").Value = 577.381
$ws.Range("N98").Value = -6896
$ws.Range("H122").Value = 1293.0416
$ws.Range("I122").Value = 920.619
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 2761.857
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -311.857
$ws.Range("N122").Value = -16600
$ws.Range("H138").Value = 2421.7058
$ws.Range("I138").Value = 2273.72
$ws.Range("J138").Value = 2564
$ws.Range("K138").Value = 6821.16
$ws.Range("L138").Value = 7692
$ws.Range("M138").Value = -1681.16
$ws.Range("N138").Value = -17972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("H102").Value = 83335060
$ws.Range("J102").Value = 333335330
$ws.Range("L102").Value = 333335330
$ws.Range("N102").Value = -333338574
$ws.Range("N44").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620
$ws.Range("H64").Value = 565.38464
$ws.Range("I64").Value = 585.7143
$ws.Range("J64").Value = 541.6667
$ws.Range("K64").Value = 585.7143
$ws.Range("L64").Value = 541.6667
$ws.Range("M64").Value = -360.7143
$ws.Range("N64").Value = -991.6667
$ws.Range("H67").Value = 565.38464
$ws.Range("I67").Value = 585.7143
$ws.Range("J67").Value = 541.6667
$ws.Range("K67").Value = 585.7143
$ws.Range("L67").Value = 541.6667
$ws.Range("M67").Value = 194.2857
$ws.Range("N67").Value = -2101.6667
$ws.Range("H82").Value = 73406.164
$ws.Range("I82").Value = 127472.8
$ws.Range("J82").Value = 34787.145
$ws.Range("K82").Value = 127472.8
$ws.Range("L82").Value = 34787.145
$ws.Range("M82").Value = -127089.8
$ws.Range("N82").Value = -35553.145
$ws.Range("H85").Value = 73406.164
$ws.Range("I85").Value = 127472.8
$ws.Range("J85").Value = 34787.145
$ws.Range("K85").Value = 127472.8
$ws.Range("L85").Value = 34787.145
$ws.Range("M85").Value = -126146.8
$ws.Range("N85").Value = -37439.145
$ws.Range("H107").Value = 3457
$ws.Range("I107").Value = 3795.1853
$ws.Range("J107").Value = 2543.9
$ws.Range("K107").Value = 3795.1853
$ws.Range("L107").Value = 2543.9
$ws.Range("M107").Value = -1875.1853
$ws.Range("N107").Value = -6383.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 14206
$ws.Range("J41").Value = 21043.334
$ws.Range("L41").Value = 21043.334
$ws.Range("N41").Value = -21899.334
$ws.Range("H50").Value = 9093.333000000001
$ws.Range("J50").Value = 9093.333000000001
$ws.Range("L50").Value = 9093.333000000001
$ws.Range("N50").Value = -10343.333
$ws.Range("H51").Value = 9721.75
$ws.Range("J51").Value = 9721.75
$ws.Range("L51").Value = 9721.75
$ws.Range("N51").Value = -11193.75
$ws.Range("H60").Value = 23408.133
$ws.Range("I60").Value = 4000
$ws.Range("J60").Value = 24794.428
$ws.Range("K60").Value = 4000
$ws.Range("L60").Value = 24794.428
$ws.Range("M60").Value = -3489
$ws.Range("N60").Value = -25816.428
$ws.Range("H61").Value = 9721.75
$ws.Range("J61").Value = 9721.75
$ws.Range("L61").Value = 9721.75
$ws.Range("N61").Value = -10417.75
$ws.Range("H68").Value = 17304
$ws.Range("J68").Value = 17304
$ws.Range("L68").Value = 17304
$ws.Range("N68").Value = -18802
$ws.Range("H71").Value = 17304
$ws.Range("J71").Value = 17304
$ws.Range("L71").Value = 51912
$ws.Range("N71").Value = -59400
$ws.Range("H109").Value = 11950
$ws.Range("J109").Value = 11950
$ws.Range("L109").Value = 11950
$ws.Range("N109").Value = -14030
$ws.Range("H132").Value = 2315.4736
$ws.Range("I132").Value = 2485.7334
$ws.Range("J132").Value = 1677
$ws.Range("K132").Value = 7457.2002
$ws.Range("L132").Value = 5031
$ws.Range("M132").Value = -4927.2002
$ws.Range("N132").Value = -10091

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 1875
$ws.Range("J19").Value = 1875
$ws.Range("L19").Value = 5625
$ws.Range("N19").Value = -5973
$ws.Range("H32").Value = 20003
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 20003
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 60009
$ws.Range("N32").Value = -60575
$ws.Range("H33").Value = 449.92307
$ws.Range("I33").Value = 431.25
$ws.Range("J33").Value = 479.8
$ws.Range("K33").Value = 2587.5
$ws.Range("L33").Value = 2878.8
$ws.Range("M33").Value = -2304.5
$ws.Range("N33").Value = -3444.8
$ws.Range("H35").Value = 3183.6428
$ws.Range("J35").Value = 3351.6155
$ws.Range("L35").Value = 10054.8465
$ws.Range("N35").Value = -10630.8465
$ws.Range("H68").Value = 1571.4166
$ws.Range("I68").Value = 4393.75
$ws.Range("J68").Value = 1006.95
$ws.Range("K68").Value = 13181.25
$ws.Range("L68").Value = 3020.85
$ws.Range("M68").Value = -12370.25
$ws.Range("N68").Value = -4642.85
$ws.Range("H71").Value = 1571.4166
$ws.Range("I71").Value = 4393.75
$ws.Range("J71").Value = 1006.95
$ws.Range("K71").Value = 39543.75
$ws.Range("L71").Value = 9062.550000000001
$ws.Range("M71").Value = -35487.75
$ws.Range("N71").Value = -17174.55
$ws.Range("H104").Value = 5200
$ws.Range("I104").Value = 600
$ws.Range("J104").Value = 7500
$ws.Range("K104").Value = 1800
$ws.Range("L104").Value = 22500
$ws.Range("M104").Value = 821
$ws.Range("N104").Value = -27742
$ws.Range("M32").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2141.3333
$ws.Range("I122").Value = 2275.75
$ws.Range("J122").Value = 1872.5
$ws.Range("K122").Value = 6827.25
$ws.Range("L122").Value = 5617.5
$ws.Range("M122").Value = -4377.25
$ws.Range("N122").Value = -10517.5
$ws.Range("H123").Value = 34425
$ws.Range("J123").Value = 34425
$ws.Range("L123").Value = 34425
$ws.Range("N123").Value = -39325
$ws.Range("H132").Value = 2034.6
$ws.Range("I132").Value = 1699.5483
$ws.Range("J132").Value = 3188.6667
$ws.Range("K132").Value = 5098.644899999999
$ws.Range("L132").Value = 9566.000100000001
$ws.Range("M132").Value = -2568.644899999999
$ws.Range("N132").Value = -14626.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7694785
$ws.Range("I7").Value = 11113134
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 11113134
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -11113022
$ws.Range("N7").Value = -3724
$ws.Range("H115").Value = 66451
$ws.Range("J115").Value = 66451
$ws.Range("L115").Value = 66451
$ws.Range("N115").Value = -68801
$ws.Range("H122").Value = 4666.6665
$ws.Range("I122").Value = 4333.3335
$ws.Range("J122").Value = 4833.3335
$ws.Range("K122").Value = 13000.0005
$ws.Range("L122").Value = 14500.0005
$ws.Range("M122").Value = -10550.0005
$ws.Range("N122").Value = -19400.0005
$ws.Range("H126").Value = 7694785
$ws.Range("I126").Value = 11113134
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 33339402
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -33336932
$ws.Range("N126").Value = -15440
$ws.Range("H132").Value = 4094.2
$ws.Range("I132").Value = 3777
$ws.Range("J132").Value = 5997.4
$ws.Range("K132").Value = 11331
$ws.Range("L132").Value = 17992.2
$ws.Range("M132").Value = -8801
$ws.Range("N132").Value = -23052.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 27588.5
$ws.Range("J109").Value = 27588.5
$ws.Range("L109").Value = 27588.5
$ws.Range("N109").Value = -30362.5
$ws.Range("H126").Value = 1005.26666
$ws.Range("I126").Value = 532.75
$ws.Range("J126").Value = 1950.3
$ws.Range("K126").Value = 1598.25
$ws.Range("L126").Value = 5850.9
$ws.Range("M126").Value = 871.75
$ws.Range("N126").Value = -10790.9
$ws.Range("H133").Value = 51307.75
$ws.Range("J133").Value = 51307.75
$ws.Range("L133").Value = 51307.75
$ws.Range("N133").Value = -61427.75
